# Statement-of-truth wording update + INCLUDEPICTURE field path normalisation.

$d = $word.ActiveDocument

# --- 1. Fix up the three INCLUDEPICTURE field instructions -----------------
# Old: /var/folders/m2/.../page1imageNNNNNNNN  \* MERGEFORMATINET
# New: C:\var\folders\m2\...\page1imageNNNNNNNN \* MERGEFORMAT
#
# The field instruction text (w:instrText) is not reachable through
# Content.Find (it is not part of the story's visible Text), so update it
# through the Fields collection instead - Field.Data rewrites the
# instruction text of the field in place.

foreach ($fld in $d.Fields) {
    $codeText = $fld.Data
    if ($codeText -and $codeText -like "*INCLUDEPICTURE*" -and $codeText -like "*/var/folders/*") {
        $newText = $codeText.Replace("/var/folders/m2/qnb2dry97b79psf_83dm0rf80000gn/T/com.microsoft.Word/WebArchiveCopyPasteTempFiles", "C:\var\folders\m2\qnb2dry97b79psf_83dm0rf80000gn\T\com.microsoft.Word\WebArchiveCopyPasteTempFiles")
        $newText = $newText.Replace("/", "\")
        $newText = $newText.Replace("MERGEFORMATINET", "MERGEFORMAT")
        $fld.Data = $newText
    }
}

# --- 2. Statement of truth wording updates ----------------------------------

$old1 = "The Claimant believes that the facts stated in the brief details of claim are true."
$new1 = "The claimant believes that the facts in this claim are true."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

$old2 = "The claimant understands that proceedings for contempt of court may be brought against anyone who makes, or causes to be made, a false statement in a document verified by a statement of truth without an honest belief in its truth."
$new2 = "The claimant understands that proceedings for contempt of court may be brought against anyone who makes, or causes to be made, a false statement in a document verified by a statement of truth without an honest belief in its truth."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
